$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.532.65'
$ws.Range("E2").Value = '  -0.79%  '

$ws.Range("D3").Value = '3.912.41'
$ws.Range("E3").Value = '  +4.20%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.16'
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.03'
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").Value = '3.912.56'
$ws.Range("E7").Value = '  +4.27%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -2.11%  '

$ws.Range("E10").Value = '  -3.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.86'
$ws.Range("E13").Value = '  -2.18%  '

$ws.Range("E14").Value = '  -1.00%  '

$ws.Range("D15").Value = '4.564.19'
$ws.Range("E15").Value = '  +4.16%  '

$ws.Range("D16").Value = '3.942.33'
$ws.Range("E16").Value = '  +4.70%  '

$ws.Range("D17").Value = '68.738.02'
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.98'
$ws.Range("E20").Value = '  -4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.19'
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '483.14'
$ws.Range("E22").Value = '  -1.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.717'
$ws.Range("E23").Value = '  -1.21%  '

$ws.Range("E24").Value = '  +13.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.35'
$ws.Range("E25").Value = '  -0.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  -1.30%  '

$ws.Range("E27").Value = '  -2.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  -0.93%  '

$ws.Range("D31").Value = '4.061.95'
$ws.Range("E31").Value = '  +4.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.84'
$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("E33").Value = '  -2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.95'
$ws.Range("E34").Value = '  +0.90%  '

$ws.Range("D35").Value = '3.854.53'
$ws.Range("E35").Value = '  +3.96%  '

$ws.Range("E36").Value = '  -0.98%  '

$ws.Range("E37").Value = '  +2.76%  '

$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("E39").Value = '  -1.04%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("E41").Value = '  -1.62%  '

$ws.Range("E42").Value = '  -3.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '433.92'
$ws.Range("E43").Value = '  +1.34%  '

$ws.Range("E44").Value = '  -0.16%  '

$ws.Range("E45").Value = '  -0.65%  '

$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").Value = '2.827.03'
$ws.Range("E48").Value = '  +0.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.83'
$ws.Range("E49").Value = '  -0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.00'
$ws.Range("E50").Value = '  +9.59%  '

$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.15'
$ws.Range("E51").Value = '  -2.60%  '

